{"js": "// Replace the two-digit-number-divided-by-one-digit-number problems\n// in the practice table with a new set of problems, preserving the\n// existing run formatting (font, size) of each cell.\nconst replacements = [\n  [\"51\u00f72=\", \"68\u00f72=\"],\n  [\"87\u00f79=\", \"17\u00f74=\"],\n  [\"16\u00f79=\", \"46\u00f79=\"],\n  [\"36\u00f78=\", \"82\u00f77=\"],\n  [\"25\u00f76=\", \"83\u00f74=\"],\n  [\"10\u00f78=\", \"54\u00f76=\"],\n  [\"18\u00f75=\", \"67\u00f72=\"],\n  [\"96\u00f75=\", \"18\u00f72=\"],\n  [\"38\u00f77=\", \"19\u00f75=\"],\n  [\"13\u00f79=\", \"65\u00f75=\"],\n  [\"23\u00f75=\", \"30\u00f73=\"],\n  [\"82\u00f73=\", \"76\u00f76=\"],\n  [\"17\u00f75=\", \"93\u00f73=\"],\n  [\"31\u00f73=\", \"69\u00f76=\"],\n  [\"52\u00f75=\", \"17\u00f74=\"],\n  [\"58\u00f75=\", \"16\u00f79=\"],\n  [\"15\u00f76=\", \"45\u00f78=\"],\n  [\"90\u00f75=\", \"50\u00f73=\"],\n  [\"64\u00f76=\", \"59\u00f74=\"],\n  [\"31\u00f77=\", \"88\u00f77=\"],\n  [\"12\u00f76=\", \"33\u00f75=\"],\n  [\"49\u00f77=\", \"39\u00f73=\"],\n  [\"44\u00f74=\", \"15\u00f72=\"],\n  [\"73\u00f78=\", \"82\u00f74=\"],\n  [\"49\u00f74=\", \"38\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items,text,font\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-number-divided-by-one-digit-number problems\n# in the practice table with a new set of problems, preserving the\n# existing run formatting (font, size) of each cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"51\u00f72=\", \"68\u00f72=\"),\n    @(\"87\u00f79=\", \"17\u00f74=\"),\n    @(\"16\u00f79=\", \"46\u00f79=\"),\n    @(\"36\u00f78=\", \"82\u00f77=\"),\n    @(\"25\u00f76=\", \"83\u00f74=\"),\n    @(\"10\u00f78=\", \"54\u00f76=\"),\n    @(\"18\u00f75=\", \"67\u00f72=\"),\n    @(\"96\u00f75=\", \"18\u00f72=\"),\n    @(\"38\u00f77=\", \"19\u00f75=\"),\n    @(\"13\u00f79=\", \"65\u00f75=\"),\n    @(\"23\u00f75=\", \"30\u00f73=\"),\n    @(\"82\u00f73=\", \"76\u00f76=\"),\n    @(\"17\u00f75=\", \"93\u00f73=\"),\n    @(\"31\u00f73=\", \"69\u00f76=\"),\n    @(\"52\u00f75=\", \"17\u00f74=\"),\n    @(\"58\u00f75=\", \"16\u00f79=\"),\n    @(\"15\u00f76=\", \"45\u00f78=\"),\n    @(\"90\u00f75=\", \"50\u00f73=\"),\n    @(\"64\u00f76=\", \"59\u00f74=\"),\n    @(\"31\u00f77=\", \"88\u00f77=\"),\n    @(\"12\u00f76=\", \"33\u00f75=\"),\n    @(\"49\u00f77=\", \"39\u00f73=\"),\n    @(\"44\u00f74=\", \"15\u00f72=\"),\n    @(\"73\u00f78=\", \"82\u00f74=\"),\n    @(\"49\u00f74=\", \"38\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # MatchCase:=True, Forward:=True, Wrap:=wdFindContinue(1), Replace:=wdReplaceAll(2)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
